$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: "{" + " " + "}" (3 runs) -> single run "{}"
# -----------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("{ }", $true, $false, $false, $false, $false, $true, 1, $false, "{}", 2) | Out-Null

# -----------------------------------------------------------------------
# Change 2: "Номер по по-" run split into "Номер по " / "по" / "-" with
# proofErr spellStart/spellEnd wrapping the inner "по".
# -----------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Номер по по-", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para2 = $rng2.Paragraphs(1).Range
$insPoint2 = $d.Range($para2.Start, $para2.Start)

$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="2E00015E" w14:textId="77777777" w:rsidR="00005576" w:rsidRPr="001D7554" w:rsidRDefault="00005576" w:rsidP="003A139E">
<w:pPr><w:tabs><w:tab w:val="left" w:pos="720"/></w:tabs><w:jc w:val="center"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>
<w:r w:rsidRPr="001D7554"><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">Номер по </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>по</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>-</w:t></w:r>
<w:r w:rsidRPr="001D7554"><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:br/><w:t>рядку</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$insPoint2.InsertXML($xml2)

# The original paragraph now sits right after the one we just inserted;
# relocate it by searching again and delete it (paragraph mark included).
$rng2b = $d.Content
$rng2b.Find.Execute("Номер по по-", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2b.Paragraphs(1).Range.Delete()

# -----------------------------------------------------------------------
# Change 3: "наимено-" / <br/> / "вание" runs gain proofErr spellStart/
# spellEnd wrapping around "наимено" and "вание" respectively, and the
# <w:br/> is split into its own run.
# -----------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("наимено-", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para3 = $rng3.Paragraphs(1).Range
$insPoint3 = $d.Range($para3.Start, $para3.Start)

$xml3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="4949D5B0" w14:textId="77777777" w:rsidR="00005576" w:rsidRPr="001D7554" w:rsidRDefault="00005576" w:rsidP="003A139E">
<w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r w:rsidRPr="001D7554"><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>наимено</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>-</w:t></w:r>
<w:r w:rsidRPr="001D7554"><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:br/></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>вание</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$insPoint3.InsertXML($xml3)

$rng3b = $d.Content
$rng3b.Find.Execute("наимено-", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng3b.Paragraphs(1).Range.Delete()

Write-Output "done"
